$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of columns B:L between row 9 and row 10
# (Mistral-7B-Instruct-v0.1 "text" and "json_object" rows were reordered)
for ($col = 2; $col -le 12; $col++) {
    $cell9 = $ws.Cells.Item(9, $col)
    $cell10 = $ws.Cells.Item(10, $col)
    $val9 = $cell9.Value()
    $val10 = $cell10.Value()
    $cell9.Value = $val10
    $cell10.Value = $val9
}
